$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a Price cell while keeping it as text (avoids Excel
# silently coercing numeric-looking strings like "287.49" into floating point
# numbers and mangling precision / trailing zeros).
function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
}

$ws.Range("D2").Value = "22.001.56"
$ws.Range("E2").Value = "  -1.43%  "

$ws.Range("D3").Value = "1.549.49"
$ws.Range("E3").Value = "  -0.98%  "

Set-TextValue "D4" "1.002"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("E5").Value = "  -0.15%  "

Set-TextValue "D6" "287.49"
$ws.Range("E6").Value = "  +0.31%  "

Set-TextValue "D7" "0.3921"
$ws.Range("E7").Value = "  +4.43%  "

$ws.Range("E8").Value = "  -2.51%  "

Set-TextValue "D9" "42.22"
$ws.Range("E9").Value = "  -7.10%  "

Set-TextValue "D10" "0.07263"
$ws.Range("E10").Value = "  -1.79%  "

Set-TextValue "D11" "1.091"
$ws.Range("E11").Value = "  -4.46%  "

Set-TextValue "D12" "1.002"
$ws.Range("E12").Value = "  -0.09%  "

$ws.Range("E13").Value = "  -7.24%  "

Set-TextValue "D14" "5.606"
$ws.Range("E14").Value = "  -3.97%  "

Set-TextValue "D15" "6.644"
$ws.Range("E15").Value = "  -2.28%  "

Set-TextValue "D16" "0.00001121"
$ws.Range("E16").Value = "  +2.68%  "

$ws.Range("D17").Value = "1.544.75"
$ws.Range("E17").Value = "  -2.13%  "

Set-TextValue "D18" "0.06579"
$ws.Range("E18").Value = "  -1.95%  "

Set-TextValue "D19" "83.66"
$ws.Range("E19").Value = "  -2.76%  "

$ws.Range("E20").Value = "  -0.03%  "

Set-TextValue "D21" "6.279"
$ws.Range("E21").Value = "  -1.03%  "

Set-TextValue "D22" "15.68"
$ws.Range("E22").Value = "  -3.46%  "

$ws.Range("E23").Value = "  -3.98%  "

$ws.Range("D24").Value = "22.002.72"
$ws.Range("E24").Value = "  -1.43%  "

$ws.Range("E25").Value = "  +2.22%  "

Set-TextValue "D26" "2.425"
$ws.Range("E26").Value = "  -3.75%  "

Set-TextValue "D27" "147.19"
$ws.Range("E27").Value = "  -1.88%  "

Set-TextValue "D28" "18.58"
$ws.Range("E28").Value = "  -4.25%  "

Set-TextValue "D29" "4.840"
$ws.Range("E29").Value = "  -0.92%  "

$ws.Range("D30").Value = "1.719.90"
$ws.Range("E30").Value = "  -1.76%  "

Set-TextValue "D31" "118.49"
$ws.Range("E31").Value = "  -4.16%  "

Set-TextValue "D32" "1.058"
$ws.Range("E32").Value = "  +0.85%  "

Set-TextValue "D33" "5.650"
$ws.Range("E33").Value = "  -4.40%  "

Set-TextValue "D34" "0.08324"
$ws.Range("E34").Value = "  +1.25%  "

Set-TextValue "D35" "9.139"
$ws.Range("E35").Value = "  -3.48%  "

Set-TextValue "D36" "1.598"
$ws.Range("E36").Value = "  -16.15%  "

Set-TextValue "D37" "0.06138"
$ws.Range("E37").Value = "  -2.32%  "

$ws.Range("E38").Value = "  -5.16%  "

Set-TextValue "D39" "5.085"
$ws.Range("E39").Value = "  -3.05%  "

Set-TextValue "D40" "1.215"
$ws.Range("E40").Value = "  -4.98%  "

Set-TextValue "D41" "0.2059"
$ws.Range("E41").Value = "  -5.43%  "

$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("E43").Value = "  -4.19%  "

Set-TextValue "D44" "0.5784"
$ws.Range("E44").Value = "  -4.71%  "

Set-TextValue "D45" "13.17"
$ws.Range("E45").Value = "  -3.81%  "

Set-TextValue "D46" "3.708"
$ws.Range("E46").Value = "  -0.83%  "

$ws.Range("E47").Value = "  -5.96%  "

$ws.Range("E48").Value = "  -4.99%  "

Set-TextValue "D49" "117.58"
$ws.Range("E49").Value = "  -4.98%  "

Set-TextValue "D50" "1.135"
$ws.Range("E50").Value = "  -3.62%  "

Set-TextValue "D51" "0.06820"
$ws.Range("E51").Value = "  -4.46%  "
